# Adds bird observation rows 105-108 to the Artfynd sheet, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 105 ---
$ws.Range("A105").Value = 131217412
$ws.Range("B105").Value = 57064
$ws.Range("D105").NumberFormat = "@"
$ws.Range("D105").Value = "NT"
$ws.Range("E105").Value = 102612
$ws.Range("F105").NumberFormat = "@"
$ws.Range("F105").Value = "Järpe"
$ws.Range("G105").NumberFormat = "@"
$ws.Range("G105").Value = "Tetrastes bonasia"
$ws.Range("H105").NumberFormat = "@"
$ws.Range("H105").Value = "(Linnaeus, 1758)"
$ws.Range("I105").NumberFormat = "@"
$ws.Range("I105").Value = "3"
$ws.Range("K105").NumberFormat = "@"
$ws.Range("K105").Value = "'"
$ws.Range("L105").NumberFormat = "@"
$ws.Range("L105").Value = "'"
$ws.Range("M105").NumberFormat = "@"
$ws.Range("M105").Value = "permanent revir"
$ws.Range("N105").NumberFormat = "@"
$ws.Range("N105").Value = "'"
$ws.Range("P105").NumberFormat = "@"
$ws.Range("P105").Value = "Garpmesliden-Römyran, Pi lm"
$ws.Range("Q105").Value = 730512
$ws.Range("R105").Value = 7281151
$ws.Range("S105").Value = 10
$ws.Range("T105").NumberFormat = "@"
$ws.Range("T105").Value = "Norrbotten"
$ws.Range("U105").NumberFormat = "@"
$ws.Range("U105").Value = "Arvidsjaur"
$ws.Range("V105").NumberFormat = "@"
$ws.Range("V105").Value = "Pite lappmark"
$ws.Range("W105").NumberFormat = "@"
$ws.Range("W105").Value = "Arvidsjaur"
$ws.Range("Y105").NumberFormat = "@"
$ws.Range("Y105").Value = "2025-12-29"
$ws.Range("AA105").NumberFormat = "@"
$ws.Range("AA105").Value = "2025-12-29"
$ws.Range("AC105").NumberFormat = "@"
$ws.Range("AC105").Value = "En tupp och två höns som satt och plockade björkknoppar. Den dels försumpade kallbäckmiljön med sin gamla kontinuitetsskog är en optimal miljö för arten. Även på vintern hittar de födda i gammelskogen."
$ws.Range("AD105").Value = $false
$ws.Range("AE105").Value = $false
$ws.Range("AG105").Value = $false
$ws.Range("AT105").NumberFormat = "@"
$ws.Range("AT105").Value = "'"
$ws.Range("AW105").NumberFormat = "@"
$ws.Range("AW105").Value = "Steve Daurer"
$ws.Range("AX105").NumberFormat = "@"
$ws.Range("AX105").Value = "Steve Daurer"
$ws.Range("AY105").NumberFormat = "@"
$ws.Range("AY105").Value = "'"

# --- Row 106 ---
$ws.Range("A106").Value = 131217375
$ws.Range("B106").Value = 58043
$ws.Range("D106").NumberFormat = "@"
$ws.Range("D106").Value = "NT"
$ws.Range("E106").Value = 103021
$ws.Range("F106").NumberFormat = "@"
$ws.Range("F106").Value = "Talltita"
$ws.Range("G106").NumberFormat = "@"
$ws.Range("G106").Value = "Poecile montanus"
$ws.Range("H106").NumberFormat = "@"
$ws.Range("H106").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I106").NumberFormat = "@"
$ws.Range("I106").Value = "2"
$ws.Range("K106").NumberFormat = "@"
$ws.Range("K106").Value = "'"
$ws.Range("L106").NumberFormat = "@"
$ws.Range("L106").Value = "'"
$ws.Range("M106").NumberFormat = "@"
$ws.Range("M106").Value = "permanent revir"
$ws.Range("N106").NumberFormat = "@"
$ws.Range("N106").Value = "'"
$ws.Range("P106").NumberFormat = "@"
$ws.Range("P106").Value = "Garpmesliden-Römyran, Pi lm"
$ws.Range("Q106").Value = 730627
$ws.Range("R106").Value = 7281017
$ws.Range("S106").Value = 10
$ws.Range("T106").NumberFormat = "@"
$ws.Range("T106").Value = "Norrbotten"
$ws.Range("U106").NumberFormat = "@"
$ws.Range("U106").Value = "Arvidsjaur"
$ws.Range("V106").NumberFormat = "@"
$ws.Range("V106").Value = "Pite lappmark"
$ws.Range("W106").NumberFormat = "@"
$ws.Range("W106").Value = "Arvidsjaur"
$ws.Range("Y106").NumberFormat = "@"
$ws.Range("Y106").Value = "2025-12-27"
$ws.Range("AA106").NumberFormat = "@"
$ws.Range("AA106").Value = "2025-12-27"
$ws.Range("AC106").NumberFormat = "@"
$ws.Range("AC106").Value = "Revirparet"
$ws.Range("AD106").Value = $false
$ws.Range("AE106").Value = $false
$ws.Range("AG106").Value = $false
$ws.Range("AT106").NumberFormat = "@"
$ws.Range("AT106").Value = "'"
$ws.Range("AW106").NumberFormat = "@"
$ws.Range("AW106").Value = "Steve Daurer"
$ws.Range("AX106").NumberFormat = "@"
$ws.Range("AX106").Value = "Steve Daurer"
$ws.Range("AY106").NumberFormat = "@"
$ws.Range("AY106").Value = "'"

# --- Row 107 ---
$ws.Range("A107").Value = 131217371
$ws.Range("B107").Value = 57988
$ws.Range("D107").NumberFormat = "@"
$ws.Range("D107").Value = "LC"
$ws.Range("E107").Value = 103031
$ws.Range("F107").NumberFormat = "@"
$ws.Range("F107").Value = "Lavskrika"
$ws.Range("G107").NumberFormat = "@"
$ws.Range("G107").Value = "Perisoreus infaustus"
$ws.Range("H107").NumberFormat = "@"
$ws.Range("H107").Value = "(Linnaeus, 1758)"
$ws.Range("I107").NumberFormat = "@"
$ws.Range("I107").Value = "4"
$ws.Range("K107").NumberFormat = "@"
$ws.Range("K107").Value = "'"
$ws.Range("L107").NumberFormat = "@"
$ws.Range("L107").Value = "'"
$ws.Range("M107").NumberFormat = "@"
$ws.Range("M107").Value = "permanent revir"
$ws.Range("N107").NumberFormat = "@"
$ws.Range("N107").Value = "'"
$ws.Range("P107").NumberFormat = "@"
$ws.Range("P107").Value = "Garpmesliden-Römyran, Pi lm"
$ws.Range("Q107").Value = 730751
$ws.Range("R107").Value = 7280826
$ws.Range("S107").Value = 10
$ws.Range("T107").NumberFormat = "@"
$ws.Range("T107").Value = "Norrbotten"
$ws.Range("U107").NumberFormat = "@"
$ws.Range("U107").Value = "Arvidsjaur"
$ws.Range("V107").NumberFormat = "@"
$ws.Range("V107").Value = "Pite lappmark"
$ws.Range("W107").NumberFormat = "@"
$ws.Range("W107").Value = "Arvidsjaur"
$ws.Range("Y107").NumberFormat = "@"
$ws.Range("Y107").Value = "2025-12-27"
$ws.Range("AA107").NumberFormat = "@"
$ws.Range("AA107").Value = "2025-12-27"
$ws.Range("AC107").NumberFormat = "@"
$ws.Range("AC107").Value = "En familjegrupp i sitt revir där det har gömt sina matförråd i de hänglavdrapperade gammelgranar så att dexklarar den bistra Lapplandvintern."
$ws.Range("AD107").Value = $false
$ws.Range("AE107").Value = $false
$ws.Range("AG107").Value = $false
$ws.Range("AT107").NumberFormat = "@"
$ws.Range("AT107").Value = "'"
$ws.Range("AW107").NumberFormat = "@"
$ws.Range("AW107").Value = "Steve Daurer"
$ws.Range("AX107").NumberFormat = "@"
$ws.Range("AX107").Value = "Steve Daurer"
$ws.Range("AY107").NumberFormat = "@"
$ws.Range("AY107").Value = "'"

# --- Row 108 ---
$ws.Range("A108").Value = 131217398
$ws.Range("B108").Value = 57884
$ws.Range("D108").NumberFormat = "@"
$ws.Range("D108").Value = "NT"
$ws.Range("E108").Value = 100109
$ws.Range("F108").NumberFormat = "@"
$ws.Range("F108").Value = "Tretåig hackspett"
$ws.Range("G108").NumberFormat = "@"
$ws.Range("G108").Value = "Picoides tridactylus"
$ws.Range("H108").NumberFormat = "@"
$ws.Range("H108").Value = "(Linnaeus, 1758)"
$ws.Range("I108").NumberFormat = "@"
$ws.Range("I108").Value = "1"
$ws.Range("K108").NumberFormat = "@"
$ws.Range("K108").Value = "adult"
$ws.Range("L108").NumberFormat = "@"
$ws.Range("L108").Value = "hane"
$ws.Range("M108").NumberFormat = "@"
$ws.Range("M108").Value = "permanent revir"
$ws.Range("N108").NumberFormat = "@"
$ws.Range("N108").Value = "'"
$ws.Range("P108").NumberFormat = "@"
$ws.Range("P108").Value = "Garpmesliden-Römyran, Pi lm"
$ws.Range("Q108").Value = 730857
$ws.Range("R108").Value = 7280781
$ws.Range("S108").Value = 10
$ws.Range("T108").NumberFormat = "@"
$ws.Range("T108").Value = "Norrbotten"
$ws.Range("U108").NumberFormat = "@"
$ws.Range("U108").Value = "Arvidsjaur"
$ws.Range("V108").NumberFormat = "@"
$ws.Range("V108").Value = "Pite lappmark"
$ws.Range("W108").NumberFormat = "@"
$ws.Range("W108").Value = "Arvidsjaur"
$ws.Range("Y108").NumberFormat = "@"
$ws.Range("Y108").Value = "2025-12-30"
$ws.Range("AA108").NumberFormat = "@"
$ws.Range("AA108").Value = "2025-12-30"
$ws.Range("AC108").NumberFormat = "@"
$ws.Range("AC108").Value = "Det noterades tidigare en revirmakerande hanne och typiska hackringar i området. Troligen samma hanne som har revir i den gamla ca 160 åriga kontinuitetsskogen."
$ws.Range("AD108").Value = $false
$ws.Range("AE108").Value = $false
$ws.Range("AG108").Value = $false
$ws.Range("AT108").NumberFormat = "@"
$ws.Range("AT108").Value = "'"
$ws.Range("AW108").NumberFormat = "@"
$ws.Range("AW108").Value = "Steve Daurer"
$ws.Range("AX108").NumberFormat = "@"
$ws.Range("AX108").Value = "Steve Daurer"
$ws.Range("AY108").NumberFormat = "@"
$ws.Range("AY108").Value = "'"

